# Generate Report for Handback
# Update the timestamp values that record when the handoff/handback
# xliff files were generated / processed.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
# Column G = "Latest HO Xliff Generate Date" for the b4821fdf... row (row 4)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-12 20:55:42"

# --- zh-cn sheet -------------------------------------------------------
# Row 4 corresponds to the b4821fdf... file.
# H = Correspond Handoff Datetime, K = Correspond Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-12 20:55:35"
$wsZhCn.Range("K4").Value = "2016-08-12 20:56:11"

# --- de-de sheet -------------------------------------------------------
# Row 4 corresponds to the b4821fdf... file.
# K = Correspond Handback DateTime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K4").Value = "2016-08-12 20:56:21"
